$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Add a new (blank) slide at the end of the deck (becomes slide 4 /
#    sldId 259), matching the new, empty ppt/slides/slide4.xml that the
#    commit introduces ("JPA Repo & setting file created").
# ---------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)   # 12 = ppLayoutBlank ("빈 화면")

# ---------------------------------------------------------------------
# 2. Re-stamp the cached "today" text of every auto-updating date field
#    ({939B0B54-38C2-4878-9BDB-D65EE6D99877}, type="datetimeFigureOut")
#    that lives on the slide master and on each of the eleven slide
#    layouts, from 2020-01-24 to 2020-01-27 (the date PowerPoint
#    recalculates such fields to whenever the deck is next saved).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2020-01-24") {
                $tr.Text = "2020-01-27"
            }
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
